$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: new log entry for "Made automatic step by step generator"
$ws.Range("A19").Value = "Made automatic step by step generator"

# B19 needs to be stored as literal text "1.5" (like several other hour
# entries in this column), not a number, while keeping the column's normal
# numeric style. Temporarily switch the cell to a text format so Excel
# stores the value verbatim, then restore the original number format.
$origFormat = $ws.Range("B19").NumberFormat
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "1.5"
$ws.Range("B19").NumberFormat = $origFormat

$ws.Range("C19").Value = 45273

$ws.Range("D19").Value = "Made a system  + UI to let the maze generate in steps automatically, so you can see the algorithm at work."
$ws.Range("D19").WrapText = $true

$ws.Range("E19").Value = "x"

$ws.Rows("19:19").RowHeight = 28.5

# Move the active selection to H18 (matches the saved view state)
$ws.Range("H18").Select() | Out-Null
